$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.447.44'
$ws.Range("E2").Value = '  +0.39%  '

$ws.Range("D3").Value = '1.908.66'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.75%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.96%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.006'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.61%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4825'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4060'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08156'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.012'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.00%  '

$ws.Range("E11").Value = '  +3.18%  '

$ws.Range("D12").Value = '1.917.90'
$ws.Range("E12").Value = '  -2.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.011'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.13%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.156'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.56%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.28'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06787'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.44%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.008'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.75%  '

$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("E19").Value = '  +0.13%  '

$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("D21").Value = '29.464.55'
$ws.Range("E21").Value = '  +0.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.622'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.73'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.73%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.185'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.80%  '

$ws.Range("D25").Value = '2.153.69'
$ws.Range("E25").Value = '  -1.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.21%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.459'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.63%  '

$ws.Range("E28").Value = '  +1.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.108'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.00%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.44'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.05%  '

$ws.Range("E31").Value = '  -3.97%  '

$ws.Range("E32").Value = '  +0.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.509'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.560'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.49%  '

$ws.Range("E35").Value = '  -2.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02269'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06104'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.175'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.24%  '

$ws.Range("E39").Value = '  +6.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5948'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.73%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.977'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.276'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.378'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.45%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07608'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.97%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.45'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.71%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5562'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.945'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '72.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.69%  '

$ws.Range("E51").Value = '  +2.48%  '

